$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the value to be stored as text (not auto-converted to a number),
    # while keeping the cell's original (default) style untouched.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.848.10"
$ws.Range("E2").Value = "  +0.39%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.236.58"
$ws.Range("E3").Value = "  +1.73%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "271.48"
$ws.Range("E5").Value = "  +4.33%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "94.31"
$ws.Range("E6").Value = "  +14.61%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.630"
$ws.Range("E7").Value = "  +1.18%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +7.96%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "46.37"
$ws.Range("E10").Value = "  +6.29%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0955"
$ws.Range("E11").Value = "  +3.92%  "

# Row 12 - Polkadot
Set-TextValue $ws.Range("D12") "8.37"
$ws.Range("E12").Value = "  +20.20%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.97%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "15.37"
$ws.Range("E14").Value = "  +7.65%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.572.08"

# Row 16 - Polygon
$ws.Range("E16").Value = "  +5.66%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.240.25"
$ws.Range("E17").Value = "  +2.46%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.807.80"
$ws.Range("E18").Value = "  +0.49%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "6.18"
$ws.Range("E20").Value = "  +4.33%  "

# Row 21 - Litecoin
Set-TextValue $ws.Range("D21") "70.87"
$ws.Range("E21").Value = "  +1.52%  "

# Row 22 - ImmutableX
$ws.Range("E22").Value = "  -4.14%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "234.71"
$ws.Range("E23").Value = "  +1.74%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "9.17"
$ws.Range("E24").Value = "  +3.41%  "

# Row 25 - Dai
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "11.47"
$ws.Range("E26").Value = "  +6.88%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  +12.08%  "

# Row 28 - WEMIXToken
$ws.Range("E28").Value = "  +6.21%  "

# Row 29 - InjectiveProtocol
Set-TextValue $ws.Range("D29") "40.33"
$ws.Range("E29").Value = "  -5.00%  "

# Row 30 - Toncoin
Set-TextValue $ws.Range("D30") "2.25"
$ws.Range("E30").Value = "  +2.53%  "

# Row 31 - Monero
$ws.Range("E31").Value = "  -0.79%  "

# Row 32 - Hedera
Set-TextValue $ws.Range("D32") "0.0917"
$ws.Range("E32").Value = "  +4.90%  "

# Row 33 - EthereumClassic
Set-TextValue $ws.Range("D33") "21.01"
$ws.Range("E33").Value = "  +2.83%  "

# Row 34 - Filecoin
Set-TextValue $ws.Range("D34") "5.50"
$ws.Range("E34").Value = "  +3.50%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +2.01%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -0.90%  "

# Row 37 - VeChain
Set-TextValue $ws.Range("D37") "0.0352"
$ws.Range("E37").Value = "  +0.22%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -3.71%  "

# Row 39 - NEARProtocol
$ws.Range("E39").Value = "  +24.35%  "

# Row 40 - Celestia
Set-TextValue $ws.Range("D40") "12.84"
$ws.Range("E40").Value = "  -1.66%  "

# Row 41 - Algorand
$ws.Range("E41").Value = "  +12.74%  "

# Row 42 - LidoDAOToken
$ws.Range("E42").Value = "  +2.48%  "

# Row 43 - MultiversX
Set-TextValue $ws.Range("D43") "63.33"
$ws.Range("E43").Value = "  -1.64%  "

# Row 44 - THORChain
$ws.Range("E44").Value = "  -1.14%  "

# Row 45 - Cronos
Set-TextValue $ws.Range("D45") "0.0997"
$ws.Range("E45").Value = "  +1.80%  "

# Row 46 - Aave
Set-TextValue $ws.Range("D46") "101.92"
$ws.Range("E46").Value = "  +1.96%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  +1.63%  "

# Row 48 - ARBITRUM
Set-TextValue $ws.Range("D48") "1.16"
$ws.Range("E48").Value = "  +4.17%  "

# Row 49 - TrustWalletToken
$ws.Range("E49").Value = "  +2.82%  "

# Row 50 - WOONetwork
Set-TextValue $ws.Range("D50") "0.450"
$ws.Range("E50").Value = "  +2.11%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.455.22"
$ws.Range("E51").Value = "  +1.72%  "
